$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.571.79"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "'1.924.81"
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "'326.27"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "'0.4821"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.08236"
$ws.Range("D10").Value = "'1.012"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'23.91"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value = "'1.919.00"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'6.122"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "'7.292"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "'0.06881"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "'1.013"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'17.65"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "'29.578.01"
$ws.Range("D22").Value = "'5.689"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Value = "'12.00"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("D24").Value = "'2.189"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'2.172.54"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "'156.33"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'6.444"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("D28").Value = "'20.05"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'2.095"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'120.91"
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("D31").Value = "'1.015"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").Value = "'0.09630"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").Value = "'5.624"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").Value = "'3.572"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "'1.381"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").Value = "'0.06404"
$ws.Range("E36").Value = "  +5.07%  "
$ws.Range("D37").Value = "'0.02296"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "'0.5964"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").Value = "'7.887"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "'0.1855"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'2.426"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").Value = "'1.281"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").Value = "'0.07548"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'0.5568"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'1.990"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").Value = "'119.88"
$ws.Range("E49").Value = "  +3.50%  "
$ws.Range("D50").Value = "'2.439"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").Value = "'72.10"
$ws.Range("E51").Value = "  -0.41%  "
